# Delete row 4 ("Miel de lavande 500 g") - all rows below shift up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(4).Delete()

# The former row 5 ("Caramiel 250 g") is now row 4; update its Image_Path
# hyperlink to point to the new caramiel image instead of the lavande one.
$ws.Hyperlinks.Item(3).Delete()
$ws.Hyperlinks.Add($ws.Range("D4"), "https://raw.githubusercontent.com/AlDenervaud/champdupuits/refs/heads/main/data/images/apiculture/caramiel_250.jpg")

# Restore the selection/active cell as recorded in the saved file.
$ws.Range("D5").Select()
